$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 7).Value = 95.63567833333333
$ws.Cells.Item(2, 8).Value = 286.907035
$ws.Cells.Item(2, 9).Value = 0.2808828217467972
$ws.Cells.Item(2, 10).Value = 0.2808828217467972
$ws.Cells.Item(2, 13).Value = 14.04766166666666
$ws.Cells.Item(2, 14).Value = 42.142985
$ws.Cells.Item(2, 15).Value = 0.3682244445652549
$ws.Cells.Item(2, 16).Value = 0.3682244445652549
$ws.Cells.Item(2, 17).Value = 1343.45765248883
$ws.Cells.Item(2, 18).Value = 12091.11887239947
$ws.Cells.Item(2, 19).Value = 0.1034279210256359
$ws.Cells.Item(2, 20).Value = 0.1034279210256359

# Row 3
$ws.Cells.Item(3, 7).Value = 95.63567833333333
$ws.Cells.Item(3, 8).Value = 286.907035
$ws.Cells.Item(3, 9).Value = 0.2808828217467972
$ws.Cells.Item(3, 10).Value = 0.2808828217467972
$ws.Cells.Item(3, 15).Value = 0.3447878236399392
$ws.Cells.Item(3, 16).Value = 0.3447878236399392
$ws.Cells.Item(3, 17).Value = 1257.949728733878
$ws.Cells.Item(3, 18).Value = 11321.54755860491
$ws.Cells.Item(3, 19).Value = 0.09684497680792317
$ws.Cells.Item(3, 20).Value = 0.09684497680792319

# Row 4
$ws.Cells.Item(4, 7).Value = 95.63567833333333
$ws.Cells.Item(4, 8).Value = 286.907035
$ws.Cells.Item(4, 9).Value = 0.2808828217467972
$ws.Cells.Item(4, 10).Value = 0.2808828217467972
$ws.Cells.Item(4, 15).Value = 0.2869877317948059
$ws.Cells.Item(4, 16).Value = 0.2869877317948059
$ws.Cells.Item(4, 17).Value = 1047.067543018095
$ws.Cells.Item(4, 18).Value = 9423.607887162852
$ws.Cells.Item(4, 19).Value = 0.08060992391323811
$ws.Cells.Item(4, 20).Value = 0.08060992391323812

# Row 5
$ws.Cells.Item(5, 9).Value = 0.392628215788982
$ws.Cells.Item(5, 10).Value = 0.392628215788982
$ws.Cells.Item(5, 13).Value = 14.04766166666666
$ws.Cells.Item(5, 14).Value = 42.142985
$ws.Cells.Item(5, 15).Value = 0.3682244445652549
$ws.Cells.Item(5, 16).Value = 0.3682244445652549
$ws.Cells.Item(5, 17).Value = 1877.933929189311
$ws.Cells.Item(5, 18).Value = 16901.4053627038
$ws.Cells.Item(5, 19).Value = 0.144575306679545
$ws.Cells.Item(5, 20).Value = 0.144575306679545

# Row 6
$ws.Cells.Item(6, 9).Value = 0.392628215788982
$ws.Cells.Item(6, 10).Value = 0.392628215788982
$ws.Cells.Item(6, 15).Value = 0.3447878236399392
$ws.Cells.Item(6, 16).Value = 0.3447878236399392
$ws.Cells.Item(6, 19).Value = 0.1353734280215155
$ws.Cells.Item(6, 20).Value = 0.1353734280215155

# Row 7
$ws.Cells.Item(7, 9).Value = 0.392628215788982
$ws.Cells.Item(7, 10).Value = 0.392628215788982
$ws.Cells.Item(7, 15).Value = 0.2869877317948059
$ws.Cells.Item(7, 16).Value = 0.2869877317948059
$ws.Cells.Item(7, 19).Value = 0.1126794810879215
$ws.Cells.Item(7, 20).Value = 0.1126794810879215

# Row 8
$ws.Cells.Item(8, 9).Value = 0.3264889624642208
$ws.Cells.Item(8, 10).Value = 0.3264889624642208
$ws.Cells.Item(8, 13).Value = 14.04766166666666
$ws.Cells.Item(8, 14).Value = 42.142985
$ws.Cells.Item(8, 15).Value = 0.3682244445652549
$ws.Cells.Item(8, 16).Value = 0.3682244445652549
$ws.Cells.Item(8, 17).Value = 1561.59103055115
$ws.Cells.Item(8, 18).Value = 14054.31927496035
$ws.Cells.Item(8, 19).Value = 0.1202212168600741
$ws.Cells.Item(8, 20).Value = 0.1202212168600741

# Row 9
$ws.Cells.Item(9, 9).Value = 0.3264889624642208
$ws.Cells.Item(9, 10).Value = 0.3264889624642208
$ws.Cells.Item(9, 15).Value = 0.3447878236399392
$ws.Cells.Item(9, 16).Value = 0.3447878236399392
$ws.Cells.Item(9, 19).Value = 0.1125694188105005
$ws.Cells.Item(9, 20).Value = 0.1125694188105005

# Row 10
$ws.Cells.Item(10, 9).Value = 0.3264889624642208
$ws.Cells.Item(10, 10).Value = 0.3264889624642208
$ws.Cells.Item(10, 15).Value = 0.2869877317948059
$ws.Cells.Item(10, 16).Value = 0.2869877317948059
$ws.Cells.Item(10, 19).Value = 0.09369832679364624
$ws.Cells.Item(10, 20).Value = 0.09369832679364626
